$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all filler pages (rows for "background", "handbell_music", "another", and
# "sample") leaving only "home" and "credit_supervised".
$ws.Range("A3:J6").EntireRow.Delete()

# The hyperlink objects do not automatically track the row deletion/shift, so clear
# the stale hyperlink collection and re-create the single hyperlink that still
# belongs on the remaining "credit_supervised" row (now row 3, column F).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/cdpeters/credit-risk-supervised-ML-sklearn") | Out-Null

# Re-apply the Hyperlink cell style that Excel normally keeps on a hyperlinked cell.
$ws.Range("F3").Style = "Hyperlink"

# Update the selected cell to match the new layout.
$ws.Range("A4").Select()
